$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - German
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

# Row 3 - Spanish
$ws.Range("F3").Value = 184
$ws.Range("G3").Value = 0.3484782485180205

# Row 11 - Indonesian
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0

# Row 12 - Finnish
$ws.Range("F12").Value = 27
$ws.Range("G12").Value = 0.1278046009656348
